# Updates the "Estado de Cuenta" worker table: the two workers' rows are
# regrouped so each worker's two periods (2209 then 2208) sit together,
# instead of alternating between workers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: SAMIRA DEL CARMEN DIAZ VARGAS - periodo 2209
$ws.Range("C16").Value = "1063144788"
$ws.Range("D16").Value = "SAMIRA DEL CARMEN DIAZ VARGAS"
$ws.Range("E16").Value = "2209"
$ws.Range("F16").Value = 16959

# Row 17: SAMIRA DEL CARMEN DIAZ VARGAS - periodo 2208
$ws.Range("C17").Value = "1063144788"
$ws.Range("D17").Value = "SAMIRA DEL CARMEN DIAZ VARGAS"
$ws.Range("E17").Value = "2208"
$ws.Range("F17").Value = 36341

# Row 18: GREGORIO JOSE LUNA FLOREZ - periodo 2209
$ws.Range("C18").Value = "1062674021"
$ws.Range("D18").Value = "GREGORIO JOSE LUNA FLOREZ"
$ws.Range("E18").Value = "2209"
$ws.Range("F18").Value = 16959

# Row 19: GREGORIO JOSE LUNA FLOREZ - periodo 2208
$ws.Range("C19").Value = "1062674021"
$ws.Range("D19").Value = "GREGORIO JOSE LUNA FLOREZ"
$ws.Range("E19").Value = "2208"
$ws.Range("F19").Value = 36341
